# Adds the CypherOutput_Message, StatOutput and StatOutput_Message worksheets
# produced by a second (stats) run of the automation tool, and refreshes the
# Cypher query text to the updated xpaths/diagnosis query that also counts
# files/samples.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the three new worksheets, in order, right after "Message".
# ---------------------------------------------------------------------------
$cypherSheet = $wb.Worksheets.Item(1)
$msgSheet    = $wb.Worksheets.Item($wb.Worksheets.Count)

$wsCypherMsg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $msgSheet)
$wsCypherMsg.Name = "CypherOutput_Message"

$wsStat = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsCypherMsg)
$wsStat.Name = "StatOutput"

$wsStatMsg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsStat)
$wsStatMsg.Name = "StatOutput_Message"

# ---------------------------------------------------------------------------
# 2. The updated Cypher query (adds OPTIONAL MATCH on file/sample + counts).
# ---------------------------------------------------------------------------
$oldQuery = $msgSheet.Range("A8").Value2
$newQuery = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN[''Lymphomatoid granulomatosis'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

# Log-message template (Neo4j connection info + cypher + output path).
$msgLines = @(
    $msgSheet.Range("A1").Value2,
    $msgSheet.Range("A2").Value2,
    $msgSheet.Range("A3").Value2,
    $msgSheet.Range("A4").Value2,
    $msgSheet.Range("A5").Value2,
    $msgSheet.Range("A6").Value2,
    $msgSheet.Range("A7").Value2,
    $oldQuery,
    $msgSheet.Range("A9").Value2,
    $msgSheet.Range("A10").Value2
)

# ---------------------------------------------------------------------------
# 3. CypherOutput_Message: exact copy of the Message log.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt $msgLines.Length; $i++) {
    $wsCypherMsg.Range("A" + ($i + 1)).Value2 = $msgLines[$i]
}

# ---------------------------------------------------------------------------
# 4. StatOutput: header row + count row (values stored as text, like source).
# ---------------------------------------------------------------------------
$wsStat.Range("A1").Value2 = "number_of_files"
$wsStat.Range("B1").Value2 = "number_of_sample"
$wsStat.Range("C1").Value2 = "number_of_cases"
$wsStat.Range("D1").Value2 = "number_of_study"

$statValues = @("0", "0", "5", "1")
for ($i = 0; $i -lt 4; $i++) {
    $cell = $wsStat.Cells.Item(2, $i + 1)
    $cell.Formula = '="' + $statValues[$i] + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# ---------------------------------------------------------------------------
# 5. StatOutput_Message: the log written twice - first run with the old
#    query, second run (after the xpath/query update) with the new query.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 10; $i++) {
    $wsStatMsg.Range("A" + ($i + 1)).Value2 = $msgLines[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $wsStatMsg.Range("A" + ($i + 11)).Value2 = $msgLines[$i]
}

$cell18 = $wsStatMsg.Range("A18")
$cell18.Formula = '="' + $newQuery.Replace('"', '""') + '"'
$cell18.Copy()
$cell18.PasteSpecial(-4163)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 6. Restore the original active sheet/selection (CypherOutput, tab 0).
# ---------------------------------------------------------------------------
$cypherSheet.Activate() | Out-Null
$cypherSheet.Range("A1").Select() | Out-Null
